$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the "*" marker in column C for the grading-element rows that were
#     missing it (or had an inconsistent style): rows 30, 31 (previously blank),
#     and 34, 35, 45 (previously no cell at all). Also re-assert it for rows
#     11, 12, 22, 32 which already held the value but used a mismatched style.
$starRows = @(11, 12, 22, 30, 31, 32, 34, 35, 45)
foreach ($r in $starRows) {
    $ws.Range("C$r").Value = "*"
}

# --- Update the saved view state of the sheet: scrolled position, zoom level
#     and the current selection.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 120
$ws.Range("B49").Select() | Out-Null
